$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: after "Install node js" add a new run ", mongodb and mongosh"
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Install node js", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$insertStart = $rng.Start
$addedText = ", mongodb and mongosh"
$rng.InsertAfter($addedText)
$insertEnd = $insertStart + $addedText.Length

# Force the newly-typed text into its own <w:r> (rather than being folded
# back into the preceding run, which would happen because the formatting
# is identical) by toggling a character attribute on just the inserted
# span and then reverting it. This mirrors the run boundary Word leaves
# behind after a real, separate editing action.
$newRun = $d.Range($insertStart, $insertEnd)
$newRun.Bold = 1
$newRun.Bold = 0

# ------------------------------------------------------------------
# Change 2: "localhost:4000/docs'" -> "localhost:4000'" and make the
# trailing quote its own run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("4000/docs", $true, $false, $false, $false, $false, $true, 1, $false, "4000", 2)

$closeQuote = [char]8217
$rng2 = $d.Content
$rng2.Find.Execute("localhost:4000" + $closeQuote, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$quoteStart = $rng2.End - 1
$quoteEnd = $rng2.End
$quoteRun = $d.Range($quoteStart, $quoteEnd)
$quoteRun.Bold = 1
$quoteRun.Bold = 0
